# Update the "Förändrad" date column (C) from 2023-09-12 (45181) to
# 2023-09-13 (45182) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
